$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prepare rows 5 and 6 by copying the formatting (styles) from row 4 ---
$ws.Range("A4:O4").Copy()
$ws.Range("A5:O5").PasteSpecial(-4122)
$ws.Range("A4:O4").Copy()
$ws.Range("A6:O6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row heights
$ws.Rows.Item(3).RowHeight = 14.9
$ws.Rows.Item(5).RowHeight = 28.35
$ws.Rows.Item(6).RowHeight = 28.35

# --- Row 3 (now "Freddy Smithers", loan count fixed, no INST1 value) ---
$ws.Range("A3").Value = 5
$ws.Range("B3").Value = 32530
$ws.Range("C3").Value = "Freddy"
$ws.Range("D3").Value = "Smithers"
$ws.Range("F3").ClearContents()
$ws.Range("H3").Value = "Cornwall"

# --- Row 4 (now "Kelly Smithson") ---
$ws.Range("A4").Value = 7
$ws.Range("B4").Value = 31423
$ws.Range("C4").Value = "Kelly"
$ws.Range("D4").Value = "Smithson"
$ws.Range("F4").Value = "Yea's Jewellers"
$ws.Range("H4").Value = "Gloucester"

# --- Row 5 (now "Tommy Smithers") ---
$ws.Range("A5").Value = 8
$ws.Range("B5").Value = 32519
$ws.Range("C5").Value = "Tommy"
$ws.Range("D5").Value = "Smithers"
$ws.Range("F5").Value = "Yea's Donuts"
$ws.Range("H5").Value = "Cornwall"
$ws.Range("I5").Value = "Ontario"
$ws.Range("J5").Value = "Canada"
$ws.Range("N5").Value = "E"

# --- Row 6 (new row, "John Smithers") ---
$ws.Range("A6").Value = 11
$ws.Range("B6").Value = 32519
$ws.Range("C6").Value = "John"
$ws.Range("D6").Value = "Smithers"
$ws.Range("F6").Value = "Yea's Donuts"
$ws.Range("H6").Value = "Cornwall"
$ws.Range("I6").Value = "Ontario"
$ws.Range("J6").Value = "Canada"
$ws.Range("N6").Value = "E"

# --- Sheet view: update dimension/selection to A7 ---
$ws.Range("A7").Select()
